# Update bitcoin_buys.xlsx after running on 2025-09-21
# Append the new day's purchase as row 49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matches the existing
# "MM/DD/YYYY" inline-string rows further up the sheet). Force a text
# number format first so Excel doesn't auto-convert the string into a
# date serial, then restore the default "Normal" style so the cell
# doesn't pick up a stray style index.
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = "09/21/2025"
$ws.Cells.Item(49, 1).Style = "Normal"

$ws.Cells.Item(49, 2).Value = 0.0004304399999999972
$ws.Cells.Item(49, 3).Value = 116160.2081590938
$ws.Cells.Item(49, 4).Value = 50
